$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text before assigning values so that numeric-looking
# strings (e.g. "1.001", "217.44") are preserved as text, matching the
# inlineStr cells in the source workbook instead of being auto-converted
# to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.472.50"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "1.676.57"
$ws.Range("E3").Value = "  +2.56%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "217.44"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").Value = "0.5317"
$ws.Range("E6").Value = "  +1.67%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +3.88%  "
$ws.Range("D9").Value = "0.06417"
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("D10").Value = "21.83"
$ws.Range("E10").Value = "  +5.74%  "
$ws.Range("D11").Value = "0.07827"
$ws.Range("E11").Value = "  +2.29%  "
$ws.Range("D12").Value = "1.681.80"
$ws.Range("E12").Value = "  +2.83%  "
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").Value = "0.5579"
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "0.0₅8340"
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("D16").Value = "65.82"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").Value = "26.509.18"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "4.760"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").Value = "194.24"
$ws.Range("E20").Value = "  +3.41%  "
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("D22").Value = "6.354"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "142.39"
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("D25").Value = "0.1288"
$ws.Range("E25").Value = "  +6.34%  "
$ws.Range("D26").Value = "7.404"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  +3.35%  "
$ws.Range("D28").Value = "1.444"
$ws.Range("E28").Value = "  +3.36%  "
$ws.Range("D29").Value = "0.06319"
$ws.Range("E29").Value = "  +6.51%  "
$ws.Range("D30").Value = "1.273"
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("D31").Value = "3.638"
$ws.Range("E31").Value = "  +5.85%  "
$ws.Range("D32").Value = "3.457"
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("D33").Value = "1.680"
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("D35").Value = "0.6212"
$ws.Range("E35").Value = "  +9.02%  "
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("D37").Value = "2.788"
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "6.175"
$ws.Range("E38").Value = "  +7.61%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01640"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("D40").Value = "1.094.97"
$ws.Range("E40").Value = "  +6.11%  "
$ws.Range("D41").Value = "0.8651"
$ws.Range("E41").Value = "  +1.63%  "
$ws.Range("D42").Value = "0.9999"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "100.53"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "1.821.99"
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("D45").Value = "58.10"
$ws.Range("E45").Value = "  +4.35%  "
$ws.Range("D46").Value = "8.208"
$ws.Range("E46").Value = "  +2.37%  "
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").Value = "0.0₈104"
$ws.Range("E48").Value = "  -5.89%  "
$ws.Range("D49").Value = "1.494"
$ws.Range("E49").Value = "  +8.09%  "
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").Value = "6.066"
$ws.Range("E51").Value = "  +2.74%  "

# Restore the default (unstyled) cell style on column D now that the
# values are stored as text, so formatting matches the original file.
$ws.Range("D2:D51").Style = "Normal"

